$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.259.58"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "1.822.59"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -1.31%  "
$ws.Range("D5").Value = "314.75"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "0.4262"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "46.02"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").Value = "0.07219"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8600"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").Value = "20.98"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "1.816.01"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "0.07107"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").Value = "87.78"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008860"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "27.293.19"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "5.131"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").Value = "2.062.20"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "2.003"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "153.17"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("D28").Value = "18.34"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "2.105"
$ws.Range("E29").Value = "  +5.32%  "
$ws.Range("D30").Value = "5.223"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "116.26"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "0.08868"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "0.7607"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "1.191"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").Value = "4.451"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.820"
$ws.Range("E36").Value = "  -6.96%  "
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "1.115"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "0.01956"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "0.05253"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.910"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Value = "7.043"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.610"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").Value = "10.54"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "106.53"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "0.4684"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").Value = "0.06397"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "1.658"
$ws.Range("E51").Value = "  -3.08%  "
